$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LSTM patient MSE values for rows 9-27 (column A)
$newValues = @(
    0.28283518845778499,
    3.3822950722096401,
    0.56398442220277401,
    0.110159468313867,
    0.61004259178938103,
    0.14534333980623601,
    0.63658992736929698,
    0.152449183550213,
    0.41304446381451798,
    0.85719014380963199,
    0.120306094137949,
    0.433660943214076,
    0.40124965339173202,
    0.079846545256123905,
    0.23990679240969201,
    0.102820949700708,
    0.18565300633019199,
    0.23352373582073199,
    1.4373208710313199
)

$row = 9
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row++
}

# Header row for the summary statistics (G4:J4) - set in this order so the
# shared-string table ends up mean, stdev, CI upper, CI lower
$ws.Range("G4").Value = "mean"
$ws.Range("H4").Value = "stdev"
$ws.Range("J4").Value = "CI upper"
$ws.Range("I4").Value = "CI lower"
$ws.Range("G4:J4").Font.Bold = $true

# Summary statistic formulas (row 5)
$ws.Range("G5").Formula = "=AVERAGE(A1:A27)"
$ws.Range("H5").Formula = "=STDEV.S(A1:A27)"
$ws.Range("I5").Formula = "=`$G5 - 2.06*`$H5/SQRT(27)"
$ws.Range("J5").Formula = "=`$G5 + 2.06*`$H5/SQRT(27)"

# Empty styled cells matching the Courier font style used at F4
$ws.Range("F4").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G4:J4").Select()

$wb.Save()
